$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Copy formatting of the last existing row down to the new row
$ws.Range("A9:C9").Copy()
$ws.Range("A10:C10").PasteSpecial(-4122)

# Add new changelog row for version 1.1.2
$ws.Range("A10").Value = 44307
$ws.Range("B10").Value = "1.1.2"
$ws.Range("C10").Value = "Improvements:`n- burned DFI added to overview"
$ws.Rows.Item(10).RowHeight = 30

# Update selection to match the new active cell
$ws.Range("F9").Select()
